# Fix IFRS financial list values (company_list sheet): correct scale/unit errors
# in rows 2-6 (2014-2018 actuals) and drop the erroneous estimate rows 7-9
# (2019E-2021E), which keep only the rank/label columns A-C.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$data = [ordered]@{
    "D" = 3377
    "E" = 61
    "F" = 61
    "G" = 58
    "H" = 40
    "I" = 41
    "J" = -1
    "K" = 2750
    "L" = 1006
    "M" = 1744
    "N" = 1744
    "O" = 0
    "P" = 154
    "Q" = 133
    "R" = -60
    "S" = -20
    "T" = 33
    "U" = 99
    "V" = 509
    "W" = 1.82
    "X" = 1.19
    "Y" = 2.36
    "Z" = 1.47
    "AA" = 57.68
    "AB" = 1067.89
    "AC" = 133
    "AD" = 21.85
    "AE" = 6785
    "AF" = 0.43
    "AG" = 50
    "AH" = 1.73
    "AI" = 31.39
    "AJ" = 30900000
}
foreach ($col in $data.Keys) {
    $ws.Range("${col}2").Value = $data[$col]
}

# Row 3
$data = [ordered]@{
    "D" = 3139
    "E" = 35
    "F" = 35
    "G" = -42
    "H" = -55
    "I" = -54
    "J" = -1
    "K" = 2565
    "L" = 903
    "M" = 1662
    "N" = 1664
    "O" = -1
    "P" = 154
    "Q" = 52
    "R" = -229
    "S" = -173
    "T" = 84
    "U" = -32
    "V" = 357
    "W" = 1.11
    "X" = -1.76
    "Y" = -3.18
    "Z" = -2.08
    "AA" = 54.3
    "AB" = 1016.28
    "AC" = -175
    "AD" = -25.18
    "AE" = 6471
    "AF" = 0.68
    "AG" = 25
    "AH" = 0.57
    "AI" = -11.86
    "AJ" = 30900000
}
foreach ($col in $data.Keys) {
    $ws.Range("${col}3").Value = $data[$col]
}

# Row 4
$data = [ordered]@{
    "D" = 2757
    "E" = 11
    "F" = 11
    "G" = 21
    "H" = 8
    "I" = 9
    "J" = -1
    "K" = 2460
    "L" = 804
    "M" = 1656
    "N" = 1656
    "O" = 0
    "P" = 154
    "Q" = 108
    "R" = 113
    "S" = -48
    "T" = 60
    "U" = 48
    "V" = 332
    "W" = 0.4
    "X" = 0.29
    "Y" = 0.54
    "Z" = 0.32
    "AA" = 48.53
    "AB" = 1014.03
    "AC" = 29
    "AD" = 111.45
    "AE" = 6442
    "AF" = 0.51
    "AG" = 25
    "AH" = 0.77
    "AI" = 71.12
    "AJ" = 30900000
}
foreach ($col in $data.Keys) {
    $ws.Range("${col}4").Value = $data[$col]
}

# Row 5
$data = [ordered]@{
    "D" = 2755
    "E" = -19
    "F" = -19
    "G" = 18
    "H" = 10
    "I" = 12
    "J" = -1
    "K" = 2514
    "L" = 860
    "M" = 1654
    "N" = 1656
    "O" = -2
    "P" = 154
    "Q" = 19
    "R" = -319
    "S" = 15
    "T" = 285
    "U" = -266
    "V" = 337
    "W" = -0.6899999999999999
    "X" = 0.38
    "Y" = 0.71
    "Z" = 0.42
    "AA" = 51.98
    "AB" = 1016.96
    "AC" = 38
    "AD" = 77.64
    "AE" = 6440
    "AF" = 0.46
    "AG" = 25
    "AH" = 0.85
    "AI" = 54.84
    "AJ" = 30900000
}
foreach ($col in $data.Keys) {
    $ws.Range("${col}5").Value = $data[$col]
}

# Row 6
$data = [ordered]@{
    "D" = 2893
    "E" = -41
    "F" = -41
    "G" = -42
    "H" = -42
    "I" = -41
    "K" = 2822
    "L" = 1230
    "M" = 1592
    "N" = 1594
    "P" = 154
    "Q" = -120
    "R" = -110
    "S" = 334
    "T" = 163
    "U" = -283
    "V" = 679
    "W" = -1.43
    "X" = -1.46
    "Y" = -2.54
    "Z" = -1.58
    "AA" = 77.28
    "AB" = 982.8200000000001
    "AC" = -134
    "AD" = -17.99
    "AE" = 6198
    "AF" = 0.39
    "AI" = 0
    "AJ" = 30900000
}
foreach ($col in $data.Keys) {
    $ws.Range("${col}6").Value = $data[$col]
}
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Rows 7-9: remove all financial figures, keep only rank / label columns A-C
$colsToClear = @("D","E","G","H","I","K","L","M","N","P","Q","R","S","T","U","W","X","Y","Z","AA","AC","AD","AE","AF","AG","AH","AI")
foreach ($row in 7..9) {
    foreach ($col in $colsToClear) {
        $ws.Range("${col}${row}").ClearContents()
    }
}
